# Refresh "cryptos" price/volume snapshot (GitHub Actions scheduled run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price ("D") column -------------------------------------------------
# Source values are plain text (e.g. "3.431.75" uses "." as a thousands
# separator), so force Text format on the whole price range first to stop
# Excel from re-parsing entries such as "1.00" or "16.90" as numbers and
# silently dropping the trailing zero / significant digits.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$priceUpdates = @{
    "D2" = "62.635.51"
    "D3" = "3.428.18"
    "D5" = "407.34"
    "D6" = "130.41"
    "D8" = "1.00"
    "D10" = "0.139"
    "D11" = "42.11"
    "D13" = "8.47"
    "D14" = "19.81"
    "D15" = "3.430.31"
    "D16" = "62.615.05"
    "D17" = "11.54"
    "D19" = "0.0000161"
    "D22" = "314.91"
    "D23" = "12.84"
    "D26" = "29.76"
    "D27" = "8.17"
    "D28" = "7.81"
    "D30" = "44.36"
    "D36" = "51.85"
    "D40" = "3.32"
    "D41" = "143.45"
    "D44" = "16.90"
    "D45" = "3.92"
    "D47" = "21.25"
    "D48" = "2.107.05"
    "D51" = "1.08"
}
foreach ($addr in $priceUpdates.Keys) {
    $ws.Range($addr).Value = $priceUpdates[$addr]
}

# Restore the default cell style now that the text is safely stored so the
# cells keep looking exactly like their untouched neighbours.
$priceRange.Style = "Normal"

# --- Volume(1h) ("E") column ---------------------------------------------
# These are already unambiguous text (leading/trailing spaces + "%"), so a
# plain value assignment is enough to keep them stored as text.
$volumeUpdates = @{
    "E2" = "  +2.22%  "
    "E3" = "  +2.52%  "
    "E4" = "  +0.13%  "
    "E5" = "  +1.25%  "
    "E6" = "  +2.48%  "
    "E7" = "  +0.80%  "
    "E8" = "  +0.01%  "
    "E9" = "  +4.81%  "
    "E10" = "  +16.07%  "
    "E11" = "  +2.34%  "
    "E12" = "  +0.28%  "
    "E13" = "  +1.87%  "
    "E14" = "  +2.37%  "
    "E15" = "  +2.79%  "
    "E16" = "  +2.41%  "
    "E17" = "  +2.25%  "
    "E18" = "  +0.92%  "
    "E19" = "  +26.16%  "
    "E20" = "  -0.75%  "
    "E21" = "  +5.30%  "
    "E22" = "  +4.76%  "
    "E23" = "  +0.16%  "
    "E24" = "  +1.75%  "
    "E25" = "  +0.43%  "
    "E26" = "  +2.48%  "
    "E27" = "  -1.64%  "
    "E28" = "  +5.02%  "
    "E29" = "  +10.23%  "
    "E30" = "  +7.70%  "
    "E31" = "  +0.67%  "
    "E32" = "  +0.87%  "
    "E33" = "  +0.46%  "
    "E34" = "  +0.00%  "
    "E35" = "  +1.50%  "
    "E36" = "  -0.51%  "
    "E37" = "  +0.28%  "
    "E38" = "  +1.85%  "
    "E39" = "  +15.01%  "
    "E40" = "  -1.74%  "
    "E41" = "  +5.27%  "
    "E42" = "  +2.53%  "
    "E43" = "  +0.72%  "
    "E44" = "  +0.94%  "
    "E45" = "  +0.92%  "
    "E46" = "  +0.31%  "
    "E47" = "  -0.49%  "
    "E48" = "  -0.64%  "
    "E49" = "  +7.02%  "
    "E50" = "  -0.95%  "
    "E51" = "  +28.44%  "
}
foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}
